$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.196.22'
$ws.Range('D3').Value = '1.686.47'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = "'216.20"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').Value = "'0.519"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = "'23.01"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.38%  '
$ws.Range('D9').Value = "'0.260"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.76%  '
$ws.Range('E10').Value = '  +0.55%  '
$ws.Range('D11').Value = "'0.0891"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').Value = '1.925.12'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = '1.681.90'
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('E14').Value = '  +2.17%  '
$ws.Range('D15').Value = "'0.557"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.05%  '
$ws.Range('D16').Value = "'66.90"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').Value = '27.199.41'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = "'236.02"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').Value = "'7.99"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.94%  '
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('E22').Value = '  +2.15%  '
$ws.Range('D23').Value = "'9.60"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.66%  '
$ws.Range('D24').Value = "'2.08"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.18%  '
$ws.Range('D25').Value = "'147.23"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').Value = "'7.37"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.58%  '
$ws.Range('D27').Value = "'16.44"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').Value = '1.542.61'
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('E35').Value = '  -1.82%  '
$ws.Range('E36').Value = '  +2.52%  '
$ws.Range('D37').Value = "'0.949"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.13%  '
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('E40').Value = '  +1.79%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'5.77"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = "'69.24"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').Value = '1.832.10'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('D47').Value = "'90.15"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('E48').Value = '  +10.07%  '
$ws.Range('D49').Value = "'1.62"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.90%  '
$ws.Range('D50').Value = "'8.28"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.20%  '
$ws.Range('E51').Value = '  -0.82%  '
